$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5620444764988075
$ws.Range("D2").Value = 0.5797676445682889

$ws.Range("C3").Value = -1.029049469922794
$ws.Range("D3").Value = 0.3146379963373782

$ws.Range("C4").Value = -0.9576317460933201
$ws.Range("D4").Value = 0.3486547911505258

$ws.Range("C5").Value = -0.5146342002074495
$ws.Range("D5").Value = 0.6119404205409649

$ws.Range("C6").Value = -1.071747440167249
$ws.Range("D6").Value = 0.2954503441901517

$ws.Range("C7").Value = -0.9350413311768538
$ws.Range("D7").Value = 0.3599177727772571

$ws.Range("C8").Value = -0.8145614527874117
$ws.Range("D8").Value = 0.424053585279232

$ws.Range("C9").Value = -0.0216649729928233
$ws.Range("D9").Value = 0.9829105050912657

$ws.Range("C10").Value = 0.2235826115837342
$ws.Range("D10").Value = 0.82514521674832

$ws.Range("C11").Value = 0.2024479331159251
$ws.Range("D11").Value = 0.8414270690090884
